# Assignment2Rubric.xlsx -- "Added story quality to the rubric"
#
# For each sheet ("Rubric" and "Grade"), insert a new requirement row
# ("Stories written correctly", worth 30/30 points) above the existing
# "At least 20 user stories" row, and rebalance the point values for
# "At least 20 user stories" (40 -> 20) and "All stories have points"
# (20 -> 10) so the rubric still totals 100 points.

$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {

    # Insert a fresh blank row above row 6 ("At least 20 user stories"),
    # pushing every row from 6 down through the old Total row down by one.
    $ws.Rows("6:6").Insert()

    # The row that used to be row 6 ("At least 20 user stories", 40/40) is
    # now row 7. Copy its formatting (and value, temporarily) back up into
    # the new row 6 so the inserted row matches the look of the other
    # requirement rows, then overwrite it with the new requirement's data.
    $ws.Range("A7:C7").Copy($ws.Range("A6:C6"))

    $ws.Range("A6").Value = "Stories written correctly"
    $ws.Range("B6").Value = 30
    $ws.Range("C6").Value = 30

    # Rebalance point totals for the two requirements whose weights changed.
    $ws.Range("B7").Value = 20
    $ws.Range("C7").Value = 20

    $ws.Range("B10").Value = 10
    $ws.Range("C10").Value = 10

    # The Total row (now row 17) should sum the new, wider range of rows,
    # including the freshly-inserted row 6.
    $ws.Range("B17").Formula = "=SUM(B6:B15)"
    $ws.Range("C17").Formula = "=SUM(C6:C15)"
}

# The "Rubric" sheet (first tab) becomes the active/selected sheet; it had
# not been previously. "Grade" (second tab) was active/selected before, so
# update its selection first, then switch to/select "Rubric" last so it
# ends up as the active tab.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Activate() | Out-Null
$ws2.Range("A6:C17").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A19").Select() | Out-Null
